$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a number
# need NumberFormat forced to Text ("@") before assignment so the stored
# value stays an exact string match (e.g. "236.36", "0.600", "55.20").

$ws.Range('D2').Value = '36.645.61'
$ws.Range('E2').Value = '  -1.99%  '
$ws.Range('D3').Value = '2.007.24'
$ws.Range('E3').Value = '  -0.14%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '236.36'
$ws.Range('E5').Value = '  -8.85%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.600'
$ws.Range('E6').Value = '  -3.12%  '
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '55.20'
$ws.Range('E8').Value = '  -3.18%  '
$ws.Range('E9').Value = '  -3.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '58.30'
$ws.Range('E10').Value = '  +3.68%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0748'
$ws.Range('E11').Value = '  -3.49%  '
$ws.Range('E12').Value = '  -3.70%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.304.66'
$ws.Range('E13').Value = '  -0.16%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.24'
$ws.Range('E14').Value = '  -0.33%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.27'
$ws.Range('E15').Value = '  -4.54%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.758'
$ws.Range('E16').Value = '  -5.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '5.12'
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').Value = '2.018.38'
$ws.Range('E18').Value = '  +0.90%  '
$ws.Range('D19').Value = '36.611.71'
$ws.Range('E19').Value = '  -1.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '67.87'
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('D21').Value = '0.0₃0806'
$ws.Range('E21').Value = '  -3.81%  '
$ws.Range('E22').Value = '  +3.79%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '222.09'
$ws.Range('E23').Value = '  -4.51%  '
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.39'
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.42'
$ws.Range('E26').Value = '  -6.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '163.23'
$ws.Range('E27').Value = '  -0.79%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.69'
$ws.Range('E28').Value = '  -3.05%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.129'
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('E30').Value = '  +3.18%  '
$ws.Range('E31').Value = '  -3.71%  '
$ws.Range('E32').Value = '  -2.47%  '
$ws.Range('E33').Value = '  -4.34%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0606'
$ws.Range('E34').Value = '  -5.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.42'
$ws.Range('E35').Value = '  +2.42%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.27'
$ws.Range('E36').Value = '  -5.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.38'
$ws.Range('E37').Value = '  +1.86%  '
$ws.Range('E38').Value = '  +0.17%  '
$ws.Range('E39').Value = '  -2.97%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.75'
$ws.Range('E40').Value = '  +6.33%  '
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('D42').Value = '1.462.99'
$ws.Range('E42').Value = '  +3.46%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0927'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0203'
$ws.Range('E44').Value = '  -4.13%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '90.47'
$ws.Range('E45').Value = '  +1.05%  '
$ws.Range('E46').Value = '  -7.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '15.38'
$ws.Range('E47').Value = '  -2.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.84'
$ws.Range('E49').Value = '  +25.51%  '
$ws.Range('E50').Value = '  -0.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.89'
$ws.Range('E51').Value = '  -2.09%  '
